$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns
$ws.Range("D2").Value = "62.981.53"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "2.563.69"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "582.31"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").Value = "142.86"
$ws.Range("E6").Value = "  -3.44%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("D9").Value = "0.105"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").Value = "5.61"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("D13").Value = "26.93"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("D14").Value = "3.024.64"
$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("D15").Value = "62.941.58"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").Value = "0.0000144"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").Value = "2.539.60"
$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("D18").Value = "10.98"
$ws.Range("E18").Value = "  -3.17%  "

$ws.Range("D19").Value = "339.31"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").Value = "4.31"
$ws.Range("E20").Value = "  -2.54%  "

$ws.Range("D21").Value = "6.59"
$ws.Range("E21").Value = "  -4.00%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  +3.58%  "

$ws.Range("D24").Value = "67.58"
$ws.Range("E24").Value = "  +1.60%  "

$ws.Range("D25").Value = "1.54"
$ws.Range("E25").Value = "  +3.27%  "

$ws.Range("D26").Value = "1.59"
$ws.Range("E26").Value = "  -2.03%  "

$ws.Range("E27").Value = "  -3.96%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "7.89"
$ws.Range("E29").Value = "  -2.99%  "

$ws.Range("D30").Value = "8.13"
$ws.Range("E30").Value = "  -3.34%  "

$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").Value = "464.08"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").Value = "0.0₃0791"
$ws.Range("E33").Value = "  -3.95%  "

$ws.Range("D34").Value = "1.65"
$ws.Range("E34").Value = "  +2.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.90"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").Value = "0.396"
$ws.Range("E37").Value = "  -2.32%  "

$ws.Range("D38").Value = "18.75"
$ws.Range("E38").Value = "  -2.07%  "

$ws.Range("D39").Value = "4.51"
$ws.Range("E39").Value = "  -0.75%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "1.68"
$ws.Range("E41").Value = "  -3.75%  "

$ws.Range("D42").Value = "40.12"
$ws.Range("E42").Value = "  +1.48%  "

$ws.Range("D43").Value = "156.49"
$ws.Range("E43").Value = "  +3.82%  "

$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  -4.33%  "

$ws.Range("D45").Value = "21.05"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("E46").Value = "  +2.65%  "

$ws.Range("D47").Value = "0.0535"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("D48").Value = "0.0958"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  -2.18%  "

$ws.Range("D50").Value = "17.98"
$ws.Range("E50").Value = "  -2.25%  "

$ws.Range("D51").Value = "11.37"
$ws.Range("E51").Value = "  -0.12%  "
